{"js": "// Replace the author name \"\u0417\u0456\u043d\u044c\u043a\u043e \u041f\u0430\u0432\u043b\u043e\" with \"\u0413\u043e\u0440\u043d\u043e\u0441\u0442\u0430\u0439 \u0411\u043e\u0433\u0434\u0430\u043d\" while\n// keeping the trailing \", \u041a\u041d-108\" group text intact, and move the\n// \"_GoBack\" bookmark (Word's \"last edit position\" marker) from the end\n// of the document to sit right after the newly-typed name - exactly\n// where Word itself drops it after a live text edit.\n\nconst body = context.document.body;\n\n// 1. Remove the old \"_GoBack\" bookmark (it currently sits in the empty\n//    paragraph at the very end of the document) so we can re-create it\n//    at the new edit location without a name clash.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Find the run of text that needs to change and swap just the name,\n//    leaving \", \u041a\u041d-108\" as-is.\nconst results = body.search(\"\u0417\u0456\u043d\u044c\u043a\u043e \u041f\u0430\u0432\u043b\u043e\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  const newRange = target.insertText(\"\u0413\u043e\u0440\u043d\u043e\u0441\u0442\u0430\u0439 \u0411\u043e\u0433\u0434\u0430\u043d\", \"Replace\");\n  await context.sync();\n\n  // 3. Drop the \"_GoBack\" bookmark right after the inserted name,\n  //    matching Word's own behaviour of marking the most recent edit.\n  const caret = newRange.getRange(\"End\");\n  caret.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Drop the stale \"_GoBack\" bookmark (currently sitting in the empty\n# paragraph at the very end of the document) so it can be re-created at\n# the actual edit location without a name clash.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Replace just the author's name, leaving the \", \u041a\u041d-108\" group intact.\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"\u0417\u0456\u043d\u044c\u043a\u043e \u041f\u0430\u0432\u043b\u043e\",     # FindText\n    $false,              # MatchCase\n    $false,              # MatchWholeWord\n    $false,              # MatchWildcards\n    $false,              # MatchSoundsLike\n    $false,              # MatchAllWordForms\n    $true,               # Forward\n    1,                   # Wrap (wdFindContinue)\n    $false,              # Format\n    \"\u0413\u043e\u0440\u043d\u043e\u0441\u0442\u0430\u0439 \u0411\u043e\u0433\u0434\u0430\u043d\",  # ReplaceWith\n    2                    # Replace (wdReplaceAll)\n)\n\n# Word re-anchors $range over the just-inserted replacement text, so\n# collapsing it to its end puts the caret right after \"\u0413\u043e\u0440\u043d\u043e\u0441\u0442\u0430\u0439 \u0411\u043e\u0433\u0434\u0430\u043d\".\n# Re-drop \"_GoBack\" there, matching Word's own behaviour of marking the\n# most recent edit location.\n$range.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $range)\n"}
